$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add scheduler in coding sections ---
# The "Biltin" item (BSL NO 2) has been dropped from the coding/scheduling
# sheet, and several items have been re-scheduled (reordered) within their
# BSL NO group. A handful of quantity cells for the rescheduled items were
# also refreshed with new figures.

# 1. Remove the "Biltin 20mg Tablet 30's" row entirely.
$ws.Rows(2).Delete()

# After the deletion above, the remaining rows shift up by one and read:
#   2 Desodin 60ml Syrup
#   3 Dinafex 180mg Tablet
#   4 Dinafex 120mg Tablet
#   5 Dinafex 60mg Tablet
#   6 Dorenta 50mg Tablet
#   7 Etorix 90mg Tablet
#   8 Etorix 120mg Tablet
#   9 Etorix 60mg Tablet - 40's
#  10 Fenobac 100ml Syrup
#  11 Flucloxin 500mg Capsule
#  12 Flucloxin 500mg Capsule - 36's
#  13 Geminox 320mg Tablet - 8's
#  14 Ketonic 30mg Injection
#  15 Ketonic 10mg Tablet
#  16 Ketonic 30mg IM/IV Injection - 4's
#  17 Kynol TR 200mg Capsule
#  18 Kynol D 25mg Tablet
#  19 Kynol TR 100mg Capsule
#  20 Naprox Plus 500mg Tablet - 30's
#  21 Oradin Plus Tablet - 40's
#  22 Osticare Tablet 24's
#  23 Rupaday Oral Solution 60ml
#  24 Sk-Mox 500mg Capsule
#  25 Zithrox 15ml Suspension
#  26 Zithrox 500mg Tablet
#  27 Zithrox 30ml Dry Suspension
#  28 Zithrox 250mg Tablet - 6's

# 2. Resequence "Dinafex 60mg Tablet" ahead of "Dinafex 180mg Tablet".
$ws.Range("A5:BF5").Copy()
$ws.Range("A3").Insert()
$ws.Rows(6).Delete()

# 3. Resequence "Etorix 120mg Tablet" ahead of "Etorix 90mg Tablet".
$ws.Range("A8:BF8").Copy()
$ws.Range("A7").Insert()
$ws.Rows(9).Delete()

# 4. Resequence "Flucloxin 500mg Capsule - 36's" ahead of "Flucloxin 500mg Capsule".
$ws.Range("A12:BF12").Copy()
$ws.Range("A11").Insert()
$ws.Rows(13).Delete()

# 5. Resequence "Ketonic 10mg Tablet" ahead of "Ketonic 30mg Injection".
$ws.Range("A15:BF15").Copy()
$ws.Range("A14").Insert()
$ws.Rows(16).Delete()

# 6. Resequence "Kynol D 25mg Tablet" ahead of "Kynol TR 200mg Capsule".
$ws.Range("A18:BF18").Copy()
$ws.Range("A17").Insert()
$ws.Rows(19).Delete()

# 7. Resequence "Zithrox 30ml Dry Suspension" ahead of "Zithrox 500mg Tablet".
$ws.Range("A27:BF27").Copy()
$ws.Range("A26").Insert()
$ws.Rows(28).Delete()

# 8. Refresh scheduled quantities for the rows that moved group position.
$ws.Range("G10").Value = 514
$ws.Range("H10").Value = 132

$ws.Range("G13").Value = 217
$ws.Range("H13").Value = 56

$ws.Range("G23").Value = 1713
$ws.Range("H23").Value = 442
$ws.Range("L23").Value = 0
$ws.Range("N23").Value = 0

$ws.Range("I24").Value = 2
$ws.Range("L24").Value = 22
$ws.Range("O24").Value = 161
$ws.Range("U24").Value = 161
$ws.Range("AM24").Value = 13
$ws.Range("AV24").Value = 148
$ws.Range("BA24").Value = 31
$ws.Range("BB24").Value = 6793
